$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.491.87"
$ws.Range("E2").Value = "  -0.60%  "
$ws.Range("D3").Value = "2.777.51"
$ws.Range("E3").Value = "  -0.08%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "352.08"
$ws.Range("E5").Value = "  -1.73%  "
$ws.Range("D6").Value = "108.25"
$ws.Range("E6").Value = "  -1.12%  "
$ws.Range("D7").Value = "0.549"
$ws.Range("E7").Value = "  -2.50%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").Value = "0.594"
$ws.Range("E9").Value = "  -0.01%  "
$ws.Range("D10").Value = "39.66"
$ws.Range("E10").Value = "  -0.73%  "
$ws.Range("E11").Value = "  +3.12%  "
$ws.Range("B12").Value = "Chainlink"
$ws.Range("C12").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D12").Value = "20.06"
$ws.Range("E12").Value = "  +3.04%  "
$ws.Range("B13").Value = "Dogecoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D13").Value = "0.0834"
$ws.Range("E13").Value = "  -1.69%  "
$ws.Range("D14").Value = "7.64"
$ws.Range("E14").Value = "  +0.74%  "
$ws.Range("D15").Value = "3.215.79"
$ws.Range("E15").Value = "  -0.03%  "
$ws.Range("D16").Value = "2.777.28"
$ws.Range("E16").Value = "  +0.84%  "
$ws.Range("E17").Value = "  -0.82%  "
$ws.Range("D18").Value = "51.503.59"
$ws.Range("E18").Value = "  -0.45%  "
$ws.Range("D19").Value = "7.63"
$ws.Range("E19").Value = "  +2.98%  "
$ws.Range("D20").Value = "3.10"
$ws.Range("E20").Value = "  -0.57%  "
$ws.Range("D21").Value = "13.12"
$ws.Range("E21").Value = "  +0.97%  "
$ws.Range("D22").Value = "0.0₃0961"
$ws.Range("E22").Value = "  -1.68%  "
$ws.Range("E23").Value = "  -0.41%  "
$ws.Range("D24").Value = "266.09"
$ws.Range("E24").Value = "  -3.02%  "
$ws.Range("D25").Value = "2.70"
$ws.Range("E25").Value = "  -1.20%  "
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("D27").Value = "26.07"
$ws.Range("E27").Value = "  -2.17%  "
$ws.Range("D28").Value = "0.161"
$ws.Range("E28").Value = "  +11.81%  "
$ws.Range("D29").Value = "10.22"
$ws.Range("E29").Value = "  +0.68%  "
$ws.Range("D30").Value = "2.22"
$ws.Range("E30").Value = "  -0.01%  "
$ws.Range("D31").Value = "36.28"
$ws.Range("E31").Value = "  +7.10%  "
$ws.Range("D32").Value = "6.22"
$ws.Range("E32").Value = "  +9.46%  "
$ws.Range("D33").Value = "51.89"
$ws.Range("E33").Value = "  +1.01%  "
$ws.Range("D34").Value = "0.0452"
$ws.Range("E34").Value = "  -2.48%  "
$ws.Range("D35").Value = "5.53"
$ws.Range("E35").Value = "  +5.53%  "
$ws.Range("E36").Value = "  -1.66%  "
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").Value = "  -0.08%  "
$ws.Range("D38").Value = "18.40"
$ws.Range("E38").Value = "  +1.71%  "
$ws.Range("D39").Value = "3.15"
$ws.Range("E39").Value = "  -2.51%  "
$ws.Range("E40").Value = "  -1.74%  "
$ws.Range("D41").Value = "2.53"
$ws.Range("E41").Value = "  -0.54%  "
$ws.Range("E42").Value = "  -0.80%  "
$ws.Range("D43").Value = "120.26"
$ws.Range("E43").Value = "  -1.26%  "
$ws.Range("D44").Value = "22.06"
$ws.Range("E44").Value = "  -0.03%  "
$ws.Range("D45").Value = "2.18"
$ws.Range("E45").Value = "  -2.68%  "
$ws.Range("D46").Value = "2.114.02"
$ws.Range("E46").Value = "  +2.20%  "
$ws.Range("D47").Value = "3.26"
$ws.Range("E47").Value = "  +0.86%  "
$ws.Range("E48").Value = "  +6.59%  "
$ws.Range("D49").Value = "5.41"
$ws.Range("E49").Value = "  -4.78%  "
$ws.Range("D50").Value = "0.901"
$ws.Range("E50").Value = "  -2.93%  "
$ws.Range("D51").Value = "1.32"
$ws.Range("E51").Value = "  +8.44%  "
